{"js": "const replacements = [\n  [\"452\u00f74=113, 0\", \"918\u00f74=229, 2\"],\n  [\"627\u00f72=313, 1\", \"681\u00f72=340, 1\"],\n  [\"433\u00f72=216, 1\", \"148\u00f75=29, 3\"],\n  [\"236\u00f75=47, 1\", \"682\u00f73=227, 1\"],\n  [\"458\u00f73=152, 2\", \"430\u00f76=71, 4\"],\n  [\"573\u00f79=63, 6\", \"314\u00f77=44, 6\"],\n  [\"797\u00f72=398, 1\", \"335\u00f77=47, 6\"],\n  [\"586\u00f76=97, 4\", \"804\u00f72=402, 0\"],\n  [\"707\u00f78=88, 3\", \"935\u00f78=116, 7\"],\n  [\"478\u00f76=79, 4\", \"309\u00f77=44, 1\"],\n  [\"210\u00f72=105, 0\", \"328\u00f79=36, 4\"],\n  [\"429\u00f73=143, 0\", \"532\u00f73=177, 1\"],\n  [\"232\u00f75=46, 2\", \"364\u00f74=91, 0\"],\n  [\"201\u00f79=22, 3\", \"688\u00f74=172, 0\"],\n  [\"795\u00f74=198, 3\", \"436\u00f79=48, 4\"],\n  [\"280\u00f77=40, 0\", \"134\u00f75=26, 4\"],\n  [\"961\u00f76=160, 1\", \"595\u00f76=99, 1\"],\n  [\"723\u00f76=120, 3\", \"379\u00f74=94, 3\"],\n  [\"795\u00f77=113, 4\", \"149\u00f79=16, 5\"],\n  [\"332\u00f79=36, 8\", \"482\u00f76=80, 2\"],\n  [\"771\u00f77=110, 1\", \"549\u00f75=109, 4\"],\n  [\"152\u00f73=50, 2\", \"158\u00f74=39, 2\"],\n  [\"789\u00f76=131, 3\", \"333\u00f79=37, 0\"],\n  [\"279\u00f74=69, 3\", \"451\u00f77=64, 3\"],\n  [\"625\u00f77=89, 2\", \"380\u00f76=63, 2\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"452\u00f74=113, 0\", \"918\u00f74=229, 2\"),\n    @(\"627\u00f72=313, 1\", \"681\u00f72=340, 1\"),\n    @(\"433\u00f72=216, 1\", \"148\u00f75=29, 3\"),\n    @(\"236\u00f75=47, 1\", \"682\u00f73=227, 1\"),\n    @(\"458\u00f73=152, 2\", \"430\u00f76=71, 4\"),\n    @(\"573\u00f79=63, 6\", \"314\u00f77=44, 6\"),\n    @(\"797\u00f72=398, 1\", \"335\u00f77=47, 6\"),\n    @(\"586\u00f76=97, 4\", \"804\u00f72=402, 0\"),\n    @(\"707\u00f78=88, 3\", \"935\u00f78=116, 7\"),\n    @(\"478\u00f76=79, 4\", \"309\u00f77=44, 1\"),\n    @(\"210\u00f72=105, 0\", \"328\u00f79=36, 4\"),\n    @(\"429\u00f73=143, 0\", \"532\u00f73=177, 1\"),\n    @(\"232\u00f75=46, 2\", \"364\u00f74=91, 0\"),\n    @(\"201\u00f79=22, 3\", \"688\u00f74=172, 0\"),\n    @(\"795\u00f74=198, 3\", \"436\u00f79=48, 4\"),\n    @(\"280\u00f77=40, 0\", \"134\u00f75=26, 4\"),\n    @(\"961\u00f76=160, 1\", \"595\u00f76=99, 1\"),\n    @(\"723\u00f76=120, 3\", \"379\u00f74=94, 3\"),\n    @(\"795\u00f77=113, 4\", \"149\u00f79=16, 5\"),\n    @(\"332\u00f79=36, 8\", \"482\u00f76=80, 2\"),\n    @(\"771\u00f77=110, 1\", \"549\u00f75=109, 4\"),\n    @(\"152\u00f73=50, 2\", \"158\u00f74=39, 2\"),\n    @(\"789\u00f76=131, 3\", \"333\u00f79=37, 0\"),\n    @(\"279\u00f74=69, 3\", \"451\u00f77=64, 3\"),\n    @(\"625\u00f77=89, 2\", \"380\u00f76=63, 2\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Find.Execute could not locate text: $oldText\"\n    }\n}\n"}
